$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the E column values for rows 18-25 (casos de prueba actualizados)
$ws.Range("E18").Value = 25
$ws.Range("E19").Value = 25
$ws.Range("E20").Value = 25
$ws.Range("E21").Value = 24
$ws.Range("E22").Value = 25
$ws.Range("E23").Value = 23
$ws.Range("E24").Value = 23
$ws.Range("E25").Value = 23

# Update the sheet view scroll position and selection
$ws.Activate()
$ws.Range("E17").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A17")
